# accuracy-per-step.xlsx edit:
# Remove the "Step" index column, leaving Technique / Intermediate Accuracy /
# (now-empty) former "Final Accuracy Without" column, and fill in the
# Intermediate Accuracy values for the remaining techniques.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Deleting column A ("Step") shifts Technique -> A, Intermediate Accuracy -> B,
# and Final Accuracy Without -> C.
$ws.Columns.Item(1).Delete()

# Fill in / correct the Intermediate Accuracy column (column B after the shift).
$ws.Range("B4").Value = 0.52
$ws.Range("B5").Value = 0.78
$ws.Range("B6").Value = 0.85
$ws.Range("B7").Value = 0.82
$ws.Range("B8").Value = 0.88
$ws.Range("B9").Value = 0.91
$ws.Range("B10").Value = 0.91
$ws.Range("B11").Value = 0.9

# The old "Final Accuracy Without" column (now column C) no longer has data.
$ws.Range("C1:C11").ClearContents()

# Match the author's final cell selection.
$ws.Range("B12").Select()
